$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Col4a5"
$ws.Range("C2").Value = "Cd93"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.022655
$ws.Range("H2").Value = 0.067965
$ws.Range("I2").Value = 0.002298240504401655
$ws.Range("J2").Value = 0.002298240504401655
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 124.722578
$ws.Range("N2").Value = 374.167734
$ws.Range("O2").Value = 0.9767311432246923
$ws.Range("P2").Value = 0.9767311432246923
$ws.Range("Q2").Value = 2.82559000459
$ws.Range("R2").Value = 25.43031004131
$ws.Range("S2").Value = 0.002244763075269522
$ws.Range("T2").Value = 0.002244763075269523

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Col4a5"
$ws.Range("C3").Value = "Cd93"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.022655
$ws.Range("H3").Value = 0.067965
$ws.Range("I3").Value = 0.002298240504401655
$ws.Range("J3").Value = 0.002298240504401655
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.134712
$ws.Range("N3").Value = 0.404136
$ws.Range("O3").Value = 0.001054960600366076
$ws.Range("P3").Value = 0.001054960600366076
$ws.Range("Q3").Value = 0.00305190036
$ws.Range("R3").Value = 0.02746710324
$ws.Range("S3").Value = 0.000002424553182309203
$ws.Range("T3").Value = 0.000002424553182309204

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Col4a5"
$ws.Range("C4").Value = "Cd93"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.022655
$ws.Range("H4").Value = 0.067965
$ws.Range("I4").Value = 0.002298240504401655
$ws.Range("J4").Value = 0.002298240504401655
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.836578333333333
$ws.Range("N4").Value = 8.509735
$ws.Range("O4").Value = 0.02221389617494163
$ws.Range("P4").Value = 0.02221389617494163
$ws.Range("Q4").Value = 0.06426268214166665
$ws.Range("R4").Value = 0.578364139275
$ws.Range("S4").Value = 0.00005105287594982384
$ws.Range("T4").Value = 0.00005105287594982385

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Col4a5"
$ws.Range("C5").Value = "Cd93"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 5.120456
$ws.Range("H5").Value = 15.361368
$ws.Range("I5").Value = 0.5194455696405421
$ws.Range("J5").Value = 0.5194455696405422
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 124.722578
$ws.Range("N5").Value = 374.167734
$ws.Range("O5").Value = 0.9767311432246923
$ws.Range("P5").Value = 0.9767311432246923
$ws.Range("Q5").Value = 638.636472855568
$ws.Range("R5").Value = 5747.728255700112
$ws.Range("S5").Value = 0.5073586650780082
$ws.Range("T5").Value = 0.5073586650780083

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Col4a5"
$ws.Range("C6").Value = "Cd93"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 5.120456
$ws.Range("H6").Value = 15.361368
$ws.Range("I6").Value = 0.5194455696405421
$ws.Range("J6").Value = 0.5194455696405422
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.134712
$ws.Range("N6").Value = 0.404136
$ws.Range("O6").Value = 0.001054960600366076
$ws.Range("P6").Value = 0.001054960600366076
$ws.Range("Q6").Value = 0.689786868672
$ws.Range("R6").Value = 6.208081818048
$ws.Range("S6").Value = 0.0005479946100054847
$ws.Range("T6").Value = 0.0005479946100054848

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Col4a5"
$ws.Range("C7").Value = "Cd93"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 5.120456
$ws.Range("H7").Value = 15.361368
$ws.Range("I7").Value = 0.5194455696405421
$ws.Range("J7").Value = 0.5194455696405422
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 2.836578333333333
$ws.Range("N7").Value = 8.509735
$ws.Range("O7").Value = 0.02221389617494163
$ws.Range("P7").Value = 0.02221389617494163
$ws.Range("Q7").Value = 14.52457454638667
$ws.Range("R7").Value = 130.72117091748
$ws.Range("S7").Value = 0.01153890995252841
$ws.Range("T7").Value = 0.01153890995252842

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Col4a5"
$ws.Range("C8").Value = "Cd93"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 4.71443
$ws.Range("H8").Value = 14.14329
$ws.Range("I8").Value = 0.4782561898550561
$ws.Range("J8").Value = 0.4782561898550561
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 124.722578
$ws.Range("N8").Value = 374.167734
$ws.Range("O8").Value = 0.9767311432246923
$ws.Range("P8").Value = 0.9767311432246923
$ws.Range("Q8").Value = 587.99586340054
$ws.Range("R8").Value = 5291.96277060486
$ws.Range("S8").Value = 0.4671277150714144
$ws.Range("T8").Value = 0.4671277150714145

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Col4a5"
$ws.Range("C9").Value = "Cd93"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 4.71443
$ws.Range("H9").Value = 14.14329
$ws.Range("I9").Value = 0.4782561898550561
$ws.Range("J9").Value = 0.4782561898550561
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.134712
$ws.Range("N9").Value = 0.404136
$ws.Range("O9").Value = 0.001054960600366076
$ws.Range("P9").Value = 0.001054960600366076
$ws.Range("Q9").Value = 0.63509029416
$ws.Range("R9").Value = 5.71581264744
$ws.Range("S9").Value = 0.000504541437178282
$ws.Range("T9").Value = 0.000504541437178282

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Col4a5"
$ws.Range("C10").Value = "Cd93"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 4.71443
$ws.Range("H10").Value = 14.14329
$ws.Range("I10").Value = 0.4782561898550561
$ws.Range("J10").Value = 0.4782561898550561
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.836578333333333
$ws.Range("N10").Value = 8.509735
$ws.Range("O10").Value = 0.02221389617494163
$ws.Range("P10").Value = 0.02221389617494163
$ws.Range("Q10").Value = 13.37284999201667
$ws.Range("R10").Value = 120.35564992815
$ws.Range("S10").Value = 0.01062393334646339
$ws.Range("T10").Value = 0.01062393334646339
